$wb = $excel.ActiveWorkbook

# Column auto-fit grows to this width (characters) once the "Status" /
# "Priority" columns contain the longer "Ready for handoff" text. The value
# is expressed relative to the ColumnWidth property (which Excel stores with
# a +5/6 character padding internally) so the saved <col width=.../> lands on
# the same figure the real Excel autofit would have produced.
$statusColWidth = 17.2159881591797 - (5/6)

# "Overview" sheet: row 3 is the 7dcc70b8... file, reporting localization status
# per target language (zh-cn column E, de-de column F) plus the latest
# handoff-xliff-generate timestamp (column G).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-03 04:16:47"
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# "zh-cn" sheet: row 3 (7dcc70b8... file) moves to "Ready for handoff" with
# machine-translation priority and an updated handoff datetime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-03 04:16:43"
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# "de-de" sheet: same row/file, same new status & priority, with its own
# (slightly later) handoff datetime.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-03 04:16:47"
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth

Write-Output "Report generated for handoff"
